# Add Betting Markets Analytics (Handicap, O/U, BTTS)
#
# - Adds three new header columns: T=Pred_Handicap, U=Pred_OU, V=Pred_BTTS
# - Fills in the actual results (K:P) for the three matches that were
#   still pending (rows 15-17: Rennes vs PSG stays pending, Atletico vs
#   Barcelona and Brentford vs Arsenal get their results)
# - Appends six new fixtures/results (rows 18-23)
# - Populates the new betting-market columns for the two most recent rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New headers T1:V1 -- copy the existing header formatting (bold, border,
#    centered) from A1 so the new cells share style index 1, then set text.
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("T1:V1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("T1").Value = "Pred_Handicap"
$ws.Range("U1").Value = "Pred_OU"
$ws.Range("V1").Value = "Pred_BTTS"

# ---------------------------------------------------------------------------
# 2. Backfill actual results for rows 15-17 (Rennes vs PSG has no actual
#    result yet, so row 15 is left alone aside from already having a
#    prediction; Atletico Madrid and Brentford vs Arsenal get filled in).
# ---------------------------------------------------------------------------

# Row 15: Rennes vs Paris S-G -> actual result now known
$ws.Range("K15").Value = "0-3"
$ws.Range("L15").Value = "Away"
$ws.Range("M15").Value = 1
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = 2
$ws.Range("P15").Value = "Sim v6.0 Correct result. Barcola hat-trick."

# Row 16: Atlético Madrid vs Barcelona -> actual result now known
$ws.Range("K16").Value = "4-0"
$ws.Range("L16").Value = "Home"
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("O16").Value = 4
$ws.Range("P16").Value = "System failure. Predicted Away/Draw, Actual Home rout."

# Row 17: Brentford vs Arsenal -> actual result now known
$ws.Range("K17").Value = "1-1"
$ws.Range("L17").Value = "Draw"
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 0
$ws.Range("O17").Value = 2
$ws.Range("P17").Value = "Arsenal dropped points. Late equalizer by Brentford."

# Helper-free approach: the Date column stores plain "YYYY-MM-DD" text (like
# the pre-existing rows), not real Excel dates. Writing that literal string
# via .Value makes Excel auto-coerce it into a date serial, so for each new
# date cell we briefly force Text format, write the string, then clear the
# formatting again (leaving the cell on the default/general style, same as
# every other text cell in the sheet) so the stored value stays literal text.

# ---------------------------------------------------------------------------
# 3. New row 18: Atlético Madrid vs Barcelona (new fixture/result entry)
# ---------------------------------------------------------------------------
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = "2026-02-12"
$ws.Range("A18").ClearFormats()
$ws.Range("B18").Value = "Atlético Madrid vs Barcelona"
$ws.Range("C18").Value = "La_Liga"
$ws.Range("D18").Value = "Atlético Madrid"
$ws.Range("E18").Value = "Barcelona"
$ws.Range("F18").Value = 31.58
$ws.Range("G18").Value = 22.51
$ws.Range("H18").Value = 45.91
$ws.Range("I18").Value = "1-1"
$ws.Range("J18").Value = "Away"
$ws.Range("K18").Value = "4-0"
$ws.Range("L18").Value = "Home"
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = 0
$ws.Range("O18").Value = 4
$ws.Range("P18").Value = "System failure. Predicted Away/Draw, Actual Home rout."

# ---------------------------------------------------------------------------
# 4. New row 19: Dortmund vs Mainz
# ---------------------------------------------------------------------------
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "2026-02-13"
$ws.Range("A19").ClearFormats()
$ws.Range("B19").Value = "Dortmund vs Mainz"
$ws.Range("C19").Value = "Bundesliga"
$ws.Range("D19").Value = "Dortmund"
$ws.Range("E19").Value = "Mainz"
$ws.Range("F19").Value = 65.5
$ws.Range("G19").Value = 20
$ws.Range("H19").Value = 14.5
$ws.Range("I19").Value = "3-1"
$ws.Range("J19").Value = "Home"
$ws.Range("K19").Value = "4-0"
$ws.Range("L19").Value = "Home"
$ws.Range("M19").Value = 1
$ws.Range("N19").Value = 0
$ws.Range("O19").Value = 2
$ws.Range("P19").Value = "Guirassy scored 2 goals. Comfortable win as predicted."

# ---------------------------------------------------------------------------
# 5. New row 20: Pisa vs Milan
# ---------------------------------------------------------------------------
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "2026-02-13"
$ws.Range("A20").ClearFormats()
$ws.Range("B20").Value = "Pisa vs Milan"
$ws.Range("C20").Value = "Serie_A"
$ws.Range("D20").Value = "Pisa"
$ws.Range("E20").Value = "Milan"
$ws.Range("F20").Value = 22.5
$ws.Range("G20").Value = 23.15
$ws.Range("H20").Value = 54.35
$ws.Range("I20").Value = "0-1"
$ws.Range("J20").Value = "Away"
$ws.Range("K20").Value = "1-2"
$ws.Range("L20").Value = "Away"
$ws.Range("M20").Value = 1
$ws.Range("N20").Value = 0
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = "Correct Result & Goal Diff. Milan won 2-1 (Modric winner)."

# ---------------------------------------------------------------------------
# 6. New row 21: Leverkusen vs St. Pauli (2026-02-14 fixture)
# ---------------------------------------------------------------------------
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "2026-02-14"
$ws.Range("A21").ClearFormats()
$ws.Range("B21").Value = "Leverkusen vs St. Pauli"
$ws.Range("C21").Value = "Bundesliga"
$ws.Range("D21").Value = "Leverkusen"
$ws.Range("E21").Value = "St. Pauli"
$ws.Range("F21").Value = 60.2
$ws.Range("G21").Value = 19.73
$ws.Range("H21").Value = 20.07
$ws.Range("I21").Value = "2-1"
$ws.Range("J21").Value = "Home"
$ws.Range("K21").Value = "4-0"
$ws.Range("L21").Value = "Home"
$ws.Range("M21").Value = 1
$ws.Range("N21").Value = 0
$ws.Range("O21").Value = 3
$ws.Range("P21").Value = "Correct Result. Leverkusen dominant win (4-0)."

# ---------------------------------------------------------------------------
# 7. New row 22: Leverkusen vs St. Pauli (2026-02-15 fixture, still pending --
#    no actual result yet) with the new betting-market predictions filled in.
# ---------------------------------------------------------------------------
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = "2026-02-15"
$ws.Range("A22").ClearFormats()
$ws.Range("B22").Value = "Leverkusen vs St. Pauli"
$ws.Range("C22").Value = "Bundesliga"
$ws.Range("D22").Value = "Leverkusen"
$ws.Range("E22").Value = "St. Pauli"
$ws.Range("F22").Value = 59.41
$ws.Range("G22").Value = 20.45
$ws.Range("H22").Value = 20.14
$ws.Range("I22").Value = "2-1"
$ws.Range("J22").Value = "Home"
$ws.Range("T22").Value = "{'Home -1.5': 36.24, 'Home -0.5': 59.41, 'Home 0.5': 79.86, 'Home 1.5': 92.55}"
$ws.Range("U22").Value = "{'Over 1.5': 82.65, 'Over 2.5': 61.46, 'Over 3.5': 39.42}"
$ws.Range("V22").Value = 59.09

# ---------------------------------------------------------------------------
# 8. New row 23: Inter vs Juventus (pending) with betting-market predictions.
# ---------------------------------------------------------------------------
$ws.Range("A23").NumberFormat = "@"
$ws.Range("A23").Value = "2026-02-15"
$ws.Range("A23").ClearFormats()
$ws.Range("B23").Value = "Inter vs Juventus"
$ws.Range("C23").Value = "Serie_A"
$ws.Range("D23").Value = "Inter"
$ws.Range("E23").Value = "Juventus"
$ws.Range("F23").Value = 45.57
$ws.Range("G23").Value = 23.95
$ws.Range("H23").Value = 30.48
$ws.Range("I23").Value = "1-1"
$ws.Range("J23").Value = "Home"
$ws.Range("T23").Value = "{'Home -1.5': 23.1, 'Home -0.5': 45.57, 'Home 0.5': 69.52000000000001, 'Home 1.5': 87.66000000000001}"
$ws.Range("U23").Value = "{'Over 1.5': 77.66, 'Over 2.5': 54.67999999999999, 'Over 3.5': 32.45}"
$ws.Range("V23").Value = 57.06
